$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update crypto price/volume/name/link data to match the latest scrape.
# Column D ("Price") values are plain numeric-looking strings (e.g. "1.001"),
# so Excel would auto-coerce them to numbers (and silently drop trailing
# zeros / collapse formatting) unless we force the cell to Text first.
# We flip NumberFormat to "@" (Text) before assigning, then reset the
# cell style back to "Normal" so no stray formatting/style diff is left
# behind - only the text content changes, matching the source diff.

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.609.96"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.62%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.804.10"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.39%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.03%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "317.94"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.32%  "

# Row 6
$ws.Range("E6").Value = "  +0.02%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5437"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.77%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3797"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.67%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07517"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.96%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "42.24"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.95%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.114"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.02%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.001"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.04%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.69"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.25%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.163"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.08%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.345"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.37%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.803.40"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.02%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "90.21"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.09%  "

# Row 18
$ws.Range("E18").Value = "  -0.68%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06549"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.06%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.47"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.14%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.954"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.68%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "28.626.55"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.62%  "

# Row 24
$ws.Range("E24").Value = "  -1.58%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.084"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.36%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "161.61"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.02%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.49"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.22%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.010.43"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.09%  "

# Row 29
$ws.Range("E29").Value = "  -3.97%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "123.01"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.97%  "

# Row 31
$ws.Range("E31").Value = "  -3.61%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1060"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.22%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.644"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.08%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.678"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.84%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.06600"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +5.11%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.2260"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.52%  "

# Row 37
$ws.Range("E37").Value = "  -0.82%  "

# Row 38
$ws.Range("B38").Value = "FraxShare"
$ws.Range("C38").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.622"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.13%  "

# Row 39
$ws.Range("B39").Value = "InternetComputer(DFINITY)"
$ws.Range("C39").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.045"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.13%  "

# Row 40
$ws.Range("B40").Value = "TheSandbox"
$ws.Range("C40").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.6230"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.64%  "

# Row 41
$ws.Range("B41").Value = "Aptos"
$ws.Range("C41").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "11.25"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.21%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.198"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.19%  "

# Row 43
$ws.Range("E43").Value = "  +4.49%  "

# Row 44
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.41"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.01%  "

# Row 45
$ws.Range("B45").Value = "PancakeSwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.695"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.19%  "

# Row 46
$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5852"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.60%  "

# Row 47
$ws.Range("B47").Value = "Quant"
$ws.Range("C47").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "127.44"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.78%  "

# Row 48
$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.962"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.72%  "

# Row 49
$ws.Range("B49").Value = "EOS"
$ws.Range("C49").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.190"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.32%  "

# Row 50
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06905"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.52%  "

# Row 51
$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "72.43"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.82%  "
